# Auto-generated edit script applying the profit-table refresh
# described in the commit diff. For each affected row, cell values
# in columns H-N are updated (changed/added/removed) to reflect the
# refreshed market-board pricing data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 4186.6665
$ws.Range("I43").Value = 3100.3333
$ws.Range("J43").Value = 4458.25
$ws.Range("K43").Value = 3100.3333
$ws.Range("L43").Value = 4458.25
$ws.Range("M43").Value = -3031.3333
$ws.Range("N43").Value = -4596.25
$ws.Range("H109").Value = 99769.42999999999
$ws.Range("J109").Value = 99769.42999999999
$ws.Range("L109").Value = 99769.42999999999
$ws.Range("N109").Value = -102543.43
$ws.Range("H113").Value = 6850
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").Value = $null
$ws.Range("H116").Value = 13485.833
$ws.Range("J116").Value = 14002.765
$ws.Range("L116").Value = 14002.765
$ws.Range("N116").Value = -20886.765
$ws.Range("H132").Value = 5981.2256
$ws.Range("I132").Value = 6367.815
$ws.Range("K132").Value = 19103.445
$ws.Range("M132").Value = -16573.445

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").Value = $null
$ws.Range("H45").Value = 3752.625
$ws.Range("I45").Value = 2000
$ws.Range("K45").Value = 2000
$ws.Range("M45").Value = -1623

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 9094145
$ws.Range("I20").Value = 14288557
$ws.Range("J20").Value = 3925
$ws.Range("K20").Value = 14288557
$ws.Range("L20").Value = 3925
$ws.Range("M20").Value = -14288310
$ws.Range("N20").Value = -4419
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").Value = $null
$ws.Range("H86").Value = 3641.9333
$ws.Range("I86").Value = 3125.75
$ws.Range("J86").Value = 3986.0557
$ws.Range("K86").Value = 3125.75
$ws.Range("L86").Value = 3986.0557
$ws.Range("M86").Value = -2002.75
$ws.Range("N86").Value = -6232.0557
$ws.Range("H89").Value = 3641.9333
$ws.Range("I89").Value = 3125.75
$ws.Range("J89").Value = 3986.0557
$ws.Range("K89").Value = 15628.75
$ws.Range("L89").Value = 19930.2785
$ws.Range("M89").Value = -10012.75
$ws.Range("N89").Value = -31162.2785
$ws.Range("H99").Value = 3879.6326
$ws.Range("I99").Value = 3612.878
$ws.Range("K99").Value = 3612.878
$ws.Range("M99").Value = -2114.878
$ws.Range("H105").Value = 3604.24
$ws.Range("I105").Value = 3348.1904
$ws.Range("J105").Value = 4948.5
$ws.Range("K105").Value = 3348.1904
$ws.Range("L105").Value = 4948.5
$ws.Range("M105").Value = -1601.1904
$ws.Range("N105").Value = -8442.5
$ws.Range("H134").Value = 3806.7
$ws.Range("I134").Value = 3451.889
$ws.Range("K134").Value = 10355.667
$ws.Range("M134").Value = -7820.667000000001
$ws.Range("H138").Value = 78837.8
$ws.Range("J138").Value = 78837.8
$ws.Range("L138").Value = 78837.8
$ws.Range("N138").Value = -89117.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 694.5
$ws.Range("I22").Value = 475.16666
$ws.Range("K22").Value = 475.16666
$ws.Range("M22").Value = -125.16666
$ws.Range("H132").Value = 1860.4
$ws.Range("I132").Value = 1860.4
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5581.200000000001
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -3051.200000000001
$ws.Range("N132").Value = $null

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 157.83333
$ws.Range("I2").Value = 51.7
$ws.Range("K2").Value = 310.2
$ws.Range("M2").Value = -197.2
$ws.Range("H7").Value = 3527.8572
$ws.Range("I7").Value = 5144.5
$ws.Range("J7").Value = 2881.2
$ws.Range("K7").Value = 15433.5
$ws.Range("L7").Value = 8643.599999999999
$ws.Range("M7").Value = -15321.5
$ws.Range("N7").Value = -8867.599999999999
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("M9").Value = $null
$ws.Range("H55").Value = 6074.5
$ws.Range("I55").Value = 916.6667
$ws.Range("J55").Value = 8285
$ws.Range("K55").Value = 2750.0001
$ws.Range("L55").Value = 24855
$ws.Range("M55").Value = -2573.0001
$ws.Range("N55").Value = -25209
$ws.Range("H107").Value = 888.3333
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").Value = $null
$ws.Range("H137").Value = 2979
$ws.Range("I137").Value = 2113.5
$ws.Range("J137").Value = 4133
$ws.Range("K137").Value = 6340.5
$ws.Range("L137").Value = 12399
$ws.Range("M137").Value = -1240.5
$ws.Range("N137").Value = -22599

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 6055555.5
$ws.Range("J11").Value = 6055555.5
$ws.Range("L11").Value = 6055555.5
$ws.Range("N11").Value = -6055833.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3604.2222
$ws.Range("I7").Value = 2705.4285
$ws.Range("J7").Value = 6750
$ws.Range("K7").Value = 2705.4285
$ws.Range("L7").Value = 6750
$ws.Range("M7").Value = -2593.4285
$ws.Range("N7").Value = -6974
$ws.Range("H38").Value = 12022
$ws.Range("J38").Value = 12022
$ws.Range("L38").Value = 12022
$ws.Range("N38").Value = -12842
$ws.Range("H82").Value = 2740.6667
$ws.Range("J82").Value = 2748.6
$ws.Range("L82").Value = 2748.6
$ws.Range("N82").Value = -3470.6
$ws.Range("H85").Value = 2740.6667
$ws.Range("J85").Value = 2748.6
$ws.Range("L85").Value = 2748.6
$ws.Range("N85").Value = -5244.6
$ws.Range("H122").Value = 4381.12
$ws.Range("I122").Value = 4092.2666
$ws.Range("K122").Value = 12276.7998
$ws.Range("M122").Value = -9826.799800000001
$ws.Range("H126").Value = 3604.2222
$ws.Range("I126").Value = 2705.4285
$ws.Range("J126").Value = 6750
$ws.Range("K126").Value = 8116.2855
$ws.Range("L126").Value = 20250
$ws.Range("M126").Value = -5646.2855
$ws.Range("N126").Value = -25190
$ws.Range("H132").Value = 41860.773
$ws.Range("I132").Value = 45185.18
$ws.Range("K132").Value = 135555.54
$ws.Range("M132").Value = -133025.54

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 1010999.7
$ws.Range("J29").Value = 16499.5
$ws.Range("L29").Value = 16499.5
$ws.Range("N29").Value = -17079.5
$ws.Range("H107").Value = 1830.6666
$ws.Range("I107").Value = 1082.6666
$ws.Range("J107").Value = 2578.6667
$ws.Range("K107").Value = 3247.9998
$ws.Range("L107").Value = 7736.000100000001
$ws.Range("M107").Value = -1327.9998
$ws.Range("N107").Value = -11576.0001
$ws.Range("H122").Value = 1998.6
$ws.Range("I122").Value = 1962.7858
$ws.Range("K122").Value = 5888.357400000001
$ws.Range("M122").Value = -3438.357400000001
$ws.Range("H136").Value = 4481.115
$ws.Range("I136").Value = 4842.35
$ws.Range("K136").Value = 14527.05
$ws.Range("M136").Value = -11977.05
